$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the recipect date (was "23.12.18" / "23 Dec 2018") to "2018-12-30".
# Go through a text formula + paste-values round trip so the result lands
# back in the cell as plain text, not an auto-recognized date serial.
$ws.Range("B2").Formula = '="2018-12-30"'
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)

$ws.Range("B3").Formula = '="2018-12-30"'
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163)

# Update the total price figures for the remaining recipects.
$ws.Range("C2").Value = 290.76
$ws.Range("C3").Value = 90.18000000000001

# Remove the cancelled sale (row 4) entirely.
$ws.Rows(4).Delete()
